{"js": "// Update the date heading paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.load(\"text\");\nawait context.sync();\n\nif (datePara.text === \"2024-05-04 Saturday\") {\n  datePara.insertText(\"2024-05-05 Sunday\", Word.InsertLocation.replace);\n}\n\n// Update every arithmetic-problem cell in the single 20x5 table, preserving\n// each cell's existing run formatting (only the text content changes).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst newValues = [\n  [\"96-49=\", \"61+22=\", \"62-4=\", \"14+74=\", \"14+61=\"],\n  [\"8+20=\", \"44+26=\", \"24+66=\", \"18+64=\", \"24+17=\"],\n  [\"99-73=\", \"59+5=\", \"98-97=\", \"48+7=\", \"26-14=\"],\n  [\"73-23=\", \"37-8=\", \"70-55=\", \"76-50=\", \"69-54=\"],\n  [\"79-73=\", \"9-4=\", \"95-18=\", \"49-7=\", \"95-42=\"],\n  [\"13+73=\", \"33+16=\", \"16+28=\", \"17+9=\", \"47-43=\"],\n  [\"90-44=\", \"42+43=\", \"32-2=\", \"46+45=\", \"4+43=\"],\n  [\"96-42=\", \"47+25=\", \"96-21=\", \"12+8=\", \"66-54=\"],\n  [\"76-48=\", \"42-9=\", \"17+56=\", \"7+65=\", \"67+20=\"],\n  [\"55+32=\", \"78-6=\", \"46+3=\", \"24+75=\", \"14+79=\"],\n  [\"43+39=\", \"77-20=\", \"18+31=\", \"17+32=\", \"27-7=\"],\n  [\"47+14=\", \"79-66=\", \"43+26=\", \"86-24=\", \"64-3=\"],\n  [\"46+25=\", \"66+29=\", \"6+33=\", \"34+54=\", \"77-69=\"],\n  [\"0+58=\", \"84-68=\", \"40+47=\", \"91-79=\", \"9+22=\"],\n  [\"92+4=\", \"54-29=\", \"96-26=\", \"56-28=\", \"54+23=\"],\n  [\"0+29=\", \"2+91=\", \"5+15=\", \"28-18=\", \"20+47=\"],\n  [\"21+49=\", \"43+49=\", \"69+1=\", \"45+30=\", \"48+12=\"],\n  [\"70+26=\", \"93-41=\", \"66-28=\", \"75-10=\", \"14+46=\"],\n  [\"0+75=\", \"74-55=\", \"80-55=\", \"58+38=\", \"10+75=\"],\n  [\"67-2=\", \"56-52=\", \"13+23=\", \"26-26=\", \"41-33=\"]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph in the document).\n$dateFind = $d.Content.Find\n$dateFind.Text = '2024-05-04 Saturday'\n$dateFind.Replacement.Text = '2024-05-05 Sunday'\n$null = $dateFind.Execute([ref]$dateFind.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$dateFind.Replacement.Text, 2)\n\n# Update every arithmetic-problem cell in the single 20x5 table in\n# document (row-major) order, preserving each cell's run formatting.\n$newValues = @(\n    @('96-49=', '61+22=', '62-4=', '14+74=', '14+61='),\n    @('8+20=', '44+26=', '24+66=', '18+64=', '24+17='),\n    @('99-73=', '59+5=', '98-97=', '48+7=', '26-14='),\n    @('73-23=', '37-8=', '70-55=', '76-50=', '69-54='),\n    @('79-73=', '9-4=', '95-18=', '49-7=', '95-42='),\n    @('13+73=', '33+16=', '16+28=', '17+9=', '47-43='),\n    @('90-44=', '42+43=', '32-2=', '46+45=', '4+43='),\n    @('96-42=', '47+25=', '96-21=', '12+8=', '66-54='),\n    @('76-48=', '42-9=', '17+56=', '7+65=', '67+20='),\n    @('55+32=', '78-6=', '46+3=', '24+75=', '14+79='),\n    @('43+39=', '77-20=', '18+31=', '17+32=', '27-7='),\n    @('47+14=', '79-66=', '43+26=', '86-24=', '64-3='),\n    @('46+25=', '66+29=', '6+33=', '34+54=', '77-69='),\n    @('0+58=', '84-68=', '40+47=', '91-79=', '9+22='),\n    @('92+4=', '54-29=', '96-26=', '56-28=', '54+23='),\n    @('0+29=', '2+91=', '5+15=', '28-18=', '20+47='),\n    @('21+49=', '43+49=', '69+1=', '45+30=', '48+12='),\n    @('70+26=', '93-41=', '66-28=', '75-10=', '14+46='),\n    @('0+75=', '74-55=', '80-55=', '58+38=', '10+75='),\n    @('67-2=', '56-52=', '13+23=', '26-26=', '41-33=')\n)\n\n$t = $d.Tables(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
